$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "week" column G: header date + weekly participation counts
$ws.Range("G1").Value = 44844
$ws.Range("G1").NumberFormat = $ws.Range("F1").NumberFormat

$ws.Range("G2").Value = 5
$ws.Range("G5").Value = 6
$ws.Range("G6").Value = 6
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 5
$ws.Range("G10").Value = 7
$ws.Range("G11").Value = 8
$ws.Range("G12").Value = 8
$ws.Range("G13").Value = 6
$ws.Range("G15").Value = 6
$ws.Range("G16").Value = 7
$ws.Range("G22").Value = 5
$ws.Range("G26").Value = 5
$ws.Range("G27").Value = 10
$ws.Range("G30").Value = 7
$ws.Range("G33").Value = 10
$ws.Range("G35").Value = 5

# Row 36 additional totals
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = 1

# Update selection to match new active cell
$ws.Range("B26").Select()
